$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 20:52"

$data = New-Object 'object[,]' 216,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 2280765
$data[0,2] = 17114
$data[0,3] = 932602
$data[0,4] = 1227090
$data[0,5] = 0
$data[0,6] = 385
$data[0,7] = 121073
$data[1,0] = 'Brasil'
$data[1,1] = 1009699
$data[1,2] = 26340
$data[1,3] = 520360
$data[1,4] = 440912
$data[1,5] = 0
$data[1,6] = 558
$data[1,7] = 48427
$data[2,0] = 'Rusia'
$data[2,1] = 569063
$data[2,2] = 7972
$data[2,3] = 324406
$data[2,4] = 236816
$data[2,5] = 0
$data[2,6] = 181
$data[2,7] = 7841
$data[3,0] = 'India'
$data[3,1] = 395812
$data[3,2] = 14721
$data[3,3] = 214206
$data[3,4] = 168636
$data[3,5] = 0
$data[3,6] = 366
$data[3,7] = 12970
$data[4,0] = 'Reino Unido'
$data[4,1] = 301815
$data[4,2] = 1346
$data[4,3] = 0
$data[4,4] = 0
$data[4,5] = 0
$data[4,6] = 173
$data[4,7] = 42461
$data[5,0] = 'España'
$data[5,1] = 292655
$data[5,2] = 307
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 27136
$data[6,0] = 'Peru'
$data[6,1] = 244388
$data[6,2] = 0
$data[6,3] = 131190
$data[6,4] = 105737
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 7461
$data[7,0] = 'Italia'
$data[7,1] = 238011
$data[7,2] = 251
$data[7,3] = 181097
$data[7,4] = 22353
$data[7,5] = 0
$data[7,6] = 47
$data[7,7] = 34561
$data[8,0] = 'Chile'
$data[8,1] = 231393
$data[8,2] = 6290
$data[8,3] = 186441
$data[8,4] = 40859
$data[8,5] = 0
$data[8,6] = 252
$data[8,7] = 4093
$data[9,0] = 'Iran'
$data[9,1] = 200262
$data[9,2] = 2615
$data[9,3] = 159192
$data[9,4] = 31678
$data[9,5] = 0
$data[9,6] = 120
$data[9,7] = 9392
$data[10,0] = 'Alemania'
$data[10,1] = 190450
$data[10,2] = 324
$data[10,3] = 174400
$data[10,4] = 7098
$data[10,5] = 0
$data[10,6] = 6
$data[10,7] = 8952
$data[11,0] = 'Turquia'
$data[11,1] = 185245
$data[11,2] = 1214
$data[11,3] = 157516
$data[11,4] = 22824
$data[11,5] = 0
$data[11,6] = 23
$data[11,7] = 4905
$data[12,0] = 'Pakistan'
$data[12,1] = 168564
$data[12,2] = 8446
$data[12,3] = 61383
$data[12,4] = 103887
$data[12,5] = 0
$data[12,6] = 201
$data[12,7] = 3294
$data[13,0] = 'Mexico'
$data[13,1] = 165455
$data[13,2] = 5662
$data[13,3] = 123095
$data[13,4] = 22613
$data[13,5] = 0
$data[13,6] = 667
$data[13,7] = 19747
$data[14,0] = 'Francia'
$data[14,1] = 159452
$data[14,2] = 811
$data[14,3] = 73887
$data[14,4] = 55948
$data[14,5] = 0
$data[14,6] = 14
$data[14,7] = 29617
$data[15,0] = 'Arabia Saudita'
$data[15,1] = 150292
$data[15,2] = 4301
$data[15,3] = 95764
$data[15,4] = 53344
$data[15,5] = 0
$data[15,6] = 45
$data[15,7] = 1184
$data[16,0] = 'Banglades'
$data[16,1] = 105535
$data[16,2] = 3243
$data[16,3] = 42945
$data[16,4] = 61202
$data[16,5] = 0
$data[16,6] = 45
$data[16,7] = 1388
$data[17,0] = 'Canada'
$data[17,1] = 100565
$data[17,2] = 345
$data[17,3] = 62961
$data[17,4] = 29258
$data[17,5] = 0
$data[17,6] = 46
$data[17,7] = 8346
$data[18,0] = 'Catar'
$data[18,1] = 85462
$data[18,2] = 1021
$data[18,3] = 65409
$data[18,4] = 19960
$data[18,5] = 0
$data[18,6] = 7
$data[18,7] = 93
$data[19,0] = 'Sudafrica'
$data[19,1] = 83890
$data[19,2] = 0
$data[19,3] = 44920
$data[19,4] = 37233
$data[19,5] = 0
$data[19,6] = 0
$data[19,7] = 1737
$data[20,0] = 'China'
$data[20,1] = 83325
$data[20,2] = 32
$data[20,3] = 78398
$data[20,4] = 293
$data[20,5] = 0
$data[20,6] = 0
$data[20,7] = 4634
$data[21,0] = 'Belgica'
$data[21,1] = 60476
$data[21,2] = 128
$data[21,3] = 16751
$data[21,4] = 34030
$data[21,5] = 0
$data[21,6] = 12
$data[21,7] = 9695
$data[22,0] = 'Colombia'
$data[22,1] = 60217
$data[22,2] = 0
$data[22,3] = 22680
$data[22,4] = 35587
$data[22,5] = 0
$data[22,6] = 0
$data[22,7] = 1950
$data[23,0] = 'Bielorrusia'
$data[23,1] = 57333
$data[23,2] = 676
$data[23,3] = 35275
$data[23,4] = 21721
$data[23,5] = 0
$data[23,6] = 6
$data[23,7] = 337
$data[24,0] = 'Suecia'
$data[24,1] = 56043
$data[24,2] = 0
$data[24,3] = 0
$data[24,4] = 0
$data[24,5] = 0
$data[24,6] = 0
$data[24,7] = 5053
$data[25,0] = 'Egipto'
$data[25,1] = 52211
$data[25,2] = 1774
$data[25,3] = 13928
$data[25,4] = 36266
$data[25,5] = 0
$data[25,6] = 79
$data[25,7] = 2017
$data[26,0] = 'Ecuador'
$data[26,1] = 49731
$data[26,2] = 634
$data[26,3] = 24446
$data[26,4] = 21129
$data[26,5] = 0
$data[26,6] = 69
$data[26,7] = 4156
$data[27,0] = 'Paises Bajos'
$data[27,1] = 49426
$data[27,2] = 107
$data[27,3] = 0
$data[27,4] = 0
$data[27,5] = 0
$data[27,6] = 3
$data[27,7] = 6081
$data[28,0] = 'Emiratos Arabes Unidos'
$data[28,1] = 44145
$data[28,2] = 393
$data[28,3] = 30996
$data[28,4] = 12849
$data[28,5] = 0
$data[28,6] = 2
$data[28,7] = 300
$data[29,0] = 'Indonesia'
$data[29,1] = 43803
$data[29,2] = 1041
$data[29,3] = 17349
$data[29,4] = 24081
$data[29,5] = 0
$data[29,6] = 34
$data[29,7] = 2373
$data[30,0] = 'Singapur'
$data[30,1] = 41615
$data[30,2] = 142
$data[30,3] = 33459
$data[30,4] = 8130
$data[30,5] = 0
$data[30,6] = 0
$data[30,7] = 26
$data[31,0] = 'Kuwait'
$data[31,1] = 38678
$data[31,2] = 604
$data[31,3] = 30190
$data[31,4] = 8175
$data[31,5] = 0
$data[31,6] = 5
$data[31,7] = 313
$data[32,0] = 'Portugal'
$data[32,1] = 38464
$data[32,2] = 375
$data[32,3] = 24477
$data[32,4] = 12460
$data[32,5] = 0
$data[32,6] = 3
$data[32,7] = 1527
$data[33,0] = 'Argentina'
$data[33,1] = 37510
$data[33,2] = 0
$data[33,3] = 11851
$data[33,4] = 24705
$data[33,5] = 0
$data[33,6] = 6
$data[33,7] = 954
$data[34,0] = 'Ucrania'
$data[34,1] = 34984
$data[34,2] = 921
$data[34,3] = 16033
$data[34,4] = 17966
$data[34,5] = 0
$data[34,6] = 19
$data[34,7] = 985
$data[35,0] = 'Polonia'
$data[35,1] = 31316
$data[35,2] = 301
$data[35,3] = 15698
$data[35,4] = 14284
$data[35,5] = 0
$data[35,6] = 18
$data[35,7] = 1334
$data[36,0] = 'Suiza'
$data[36,1] = 31217
$data[36,2] = 17
$data[36,3] = 28900
$data[36,4] = 361
$data[36,5] = 0
$data[36,6] = 0
$data[36,7] = 1956
$data[37,0] = 'Filipinas'
$data[37,1] = 28459
$data[37,2] = 660
$data[37,3] = 7378
$data[37,4] = 19951
$data[37,5] = 0
$data[37,6] = 14
$data[37,7] = 1130
$data[38,0] = 'Afganistan'
$data[38,1] = 27878
$data[38,2] = 346
$data[38,3] = 7962
$data[38,4] = 19368
$data[38,5] = 0
$data[38,6] = 2
$data[38,7] = 548
$data[39,0] = 'Oman'
$data[39,1] = 27670
$data[39,2] = 852
$data[39,3] = 13974
$data[39,4] = 13571
$data[39,5] = 0
$data[39,6] = 6
$data[39,7] = 125
$data[40,0] = 'Irak'
$data[40,1] = 27352
$data[40,2] = 1635
$data[40,3] = 12205
$data[40,4] = 14222
$data[40,5] = 0
$data[40,6] = 69
$data[40,7] = 925
$data[41,0] = 'Irlanda'
$data[41,1] = 25368
$data[41,2] = 13
$data[41,3] = 22698
$data[41,4] = 956
$data[41,5] = 0
$data[41,6] = 0
$data[41,7] = 1714
$data[42,0] = 'Republica Dominicana'
$data[42,1] = 25068
$data[42,2] = 423
$data[42,3] = 14605
$data[42,4] = 9816
$data[42,5] = 0
$data[42,6] = 12
$data[42,7] = 647
$data[43,0] = 'Rumania'
$data[43,1] = 23400
$data[43,2] = 320
$data[43,3] = 16555
$data[43,4] = 5361
$data[43,5] = 0
$data[43,6] = 11
$data[43,7] = 1484
$data[44,0] = 'Panama'
$data[44,1] = 23351
$data[44,2] = 0
$data[44,3] = 13782
$data[44,4] = 9094
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 475
$data[45,0] = 'Bolivia'
$data[45,1] = 21499
$data[45,2] = 814
$data[45,3] = 4002
$data[45,4] = 16800
$data[45,5] = 0
$data[45,6] = 18
$data[45,7] = 697
$data[46,0] = 'Barein'
$data[46,1] = 20430
$data[46,2] = 0
$data[46,3] = 14696
$data[46,4] = 5677
$data[46,5] = 0
$data[46,6] = 2
$data[46,7] = 57
$data[47,0] = 'Israel'
$data[47,1] = 20339
$data[47,2] = 303
$data[47,3] = 15586
$data[47,4] = 4449
$data[47,5] = 0
$data[47,6] = 1
$data[47,7] = 304
$data[48,0] = 'Armenia'
$data[48,1] = 19157
$data[48,2] = 459
$data[48,3] = 8266
$data[48,4] = 10572
$data[48,5] = 0
$data[48,6] = 10
$data[48,7] = 319
$data[49,0] = 'Nigeria'
$data[49,1] = 18480
$data[49,2] = 0
$data[49,3] = 6307
$data[49,4] = 11698
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 475
$data[50,0] = 'Japon'
$data[50,1] = 17668
$data[50,2] = 0
$data[50,3] = 15930
$data[50,4] = 803
$data[50,5] = 0
$data[50,6] = 0
$data[50,7] = 935
$data[51,0] = 'Austria'
$data[51,1] = 17271
$data[51,2] = 48
$data[51,3] = 16141
$data[51,4] = 442
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 688
$data[52,0] = 'Kazajistan'
$data[52,1] = 16351
$data[52,2] = 474
$data[52,3] = 10195
$data[52,4] = 6043
$data[52,5] = 0
$data[52,6] = 13
$data[52,7] = 113
$data[53,0] = 'Moldavia'
$data[53,1] = 13556
$data[53,2] = 450
$data[53,3] = 7525
$data[53,4] = 5581
$data[53,5] = 0
$data[53,6] = 6
$data[53,7] = 450
$data[54,0] = 'Ghana'
$data[54,1] = 13203
$data[54,2] = 274
$data[54,3] = 4548
$data[54,4] = 8585
$data[54,5] = 0
$data[54,6] = 4
$data[54,7] = 70
$data[55,0] = 'Serbia'
$data[55,1] = 12709
$data[55,2] = 93
$data[55,3] = 11511
$data[55,4] = 939
$data[55,5] = 0
$data[55,6] = 1
$data[55,7] = 259
$data[56,0] = 'Dinamarca'
$data[56,1] = 12391
$data[56,2] = 47
$data[56,3] = 11282
$data[56,4] = 509
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 600
$data[57,0] = 'Corea del Sur'
$data[57,1] = 12306
$data[57,2] = 49
$data[57,3] = 10835
$data[57,4] = 1191
$data[57,5] = 0
$data[57,6] = 0
$data[57,7] = 280
$data[58,0] = 'Guatemala'
$data[58,1] = 11868
$data[58,2] = 617
$data[58,3] = 2290
$data[58,4] = 9129
$data[58,5] = 0
$data[58,6] = 17
$data[58,7] = 449
$data[59,0] = 'Azerbaiyan'
$data[59,1] = 11767
$data[59,2] = 438
$data[59,3] = 6325
$data[59,4] = 5299
$data[59,5] = 0
$data[59,6] = 4
$data[59,7] = 143
$data[60,0] = 'Argelia'
$data[60,1] = 11504
$data[60,2] = 119
$data[60,3] = 8196
$data[60,4] = 2483
$data[60,5] = 0
$data[60,6] = 14
$data[60,7] = 825
$data[61,0] = 'Honduras'
$data[61,1] = 10739
$data[61,2] = 440
$data[61,3] = 1179
$data[61,4] = 9217
$data[61,5] = 0
$data[61,6] = 7
$data[61,7] = 343
$data[62,0] = 'Camerun'
$data[62,1] = 10638
$data[62,2] = 774
$data[62,3] = 7548
$data[62,4] = 2808
$data[62,5] = 0
$data[62,6] = 6
$data[62,7] = 282
$data[63,0] = 'Chequia'
$data[63,1] = 10361
$data[63,2] = 81
$data[63,3] = 7472
$data[63,4] = 2554
$data[63,5] = 0
$data[63,6] = 1
$data[63,7] = 335
$data[64,0] = 'Marruecos'
$data[64,1] = 9613
$data[64,2] = 539
$data[64,3] = 8117
$data[64,4] = 1283
$data[64,5] = 0
$data[64,6] = 0
$data[64,7] = 213
$data[65,0] = 'Noruega'
$data[65,1] = 8721
$data[65,2] = 13
$data[65,3] = 8138
$data[65,4] = 339
$data[65,5] = 0
$data[65,6] = 0
$data[65,7] = 244
$data[66,0] = 'Malasia'
$data[66,1] = 8535
$data[66,2] = 6
$data[66,3] = 8070
$data[66,4] = 344
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 121
$data[67,0] = 'Sudan'
$data[67,1] = 8316
$data[67,2] = 296
$data[67,3] = 3086
$data[67,4] = 4724
$data[67,5] = 0
$data[67,6] = 19
$data[67,7] = 506
$data[68,0] = 'Nepal'
$data[68,1] = 8274
$data[68,2] = 426
$data[68,3] = 1402
$data[68,4] = 6850
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 22
$data[69,0] = 'Australia'
$data[69,1] = 7409
$data[69,2] = 18
$data[69,3] = 6878
$data[69,4] = 429
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 102
$data[70,0] = 'Finlandia'
$data[70,1] = 7133
$data[70,2] = 14
$data[70,3] = 6200
$data[70,4] = 607
$data[70,5] = 0
$data[70,6] = 0
$data[70,7] = 326
$data[71,0] = 'Costa de Marfil'
$data[71,1] = 6444
$data[71,2] = 0
$data[71,3] = 2863
$data[71,4] = 3532
$data[71,5] = 0
$data[71,6] = 0
$data[71,7] = 49
$data[72,0] = 'Uzbekistan'
$data[72,1] = 5920
$data[72,2] = 153
$data[72,3] = 4273
$data[72,4] = 1628
$data[72,5] = 0
$data[72,6] = 0
$data[72,7] = 19
$data[73,0] = 'Senegal'
$data[73,1] = 5639
$data[73,2] = 164
$data[73,3] = 3788
$data[73,4] = 1772
$data[73,5] = 0
$data[73,6] = 3
$data[73,7] = 79
$data[74,0] = 'Consejo Danes para los Refugiados'
$data[74,1] = 5477
$data[74,2] = 194
$data[74,3] = 719
$data[74,4] = 4636
$data[74,5] = 0
$data[74,6] = 5
$data[74,7] = 122
$data[75,0] = 'Tayikistan'
$data[75,1] = 5338
$data[75,2] = 59
$data[75,3] = 3830
$data[75,4] = 1457
$data[75,5] = 0
$data[75,6] = 0
$data[75,7] = 51
$data[76,0] = 'Haiti'
$data[76,1] = 4916
$data[76,2] = 228
$data[76,3] = 24
$data[76,4] = 4808
$data[76,5] = 0
$data[76,6] = 2
$data[76,7] = 84
$data[77,0] = 'Guinea'
$data[77,1] = 4841
$data[77,2] = 0
$data[77,3] = 3467
$data[77,4] = 1348
$data[77,5] = 0
$data[77,6] = 0
$data[77,7] = 26
$data[78,0] = 'Republica de Macedonia'
$data[78,1] = 4820
$data[78,2] = 156
$data[78,3] = 1863
$data[78,4] = 2735
$data[78,5] = 0
$data[78,6] = 6
$data[78,7] = 222
$data[79,0] = 'Republica de Yibuti'
$data[79,1] = 4565
$data[79,2] = 8
$data[79,3] = 3565
$data[79,4] = 955
$data[79,5] = 0
$data[79,6] = 2
$data[79,7] = 45
$data[80,0] = 'Kenia'
$data[80,1] = 4374
$data[80,2] = 117
$data[80,3] = 1550
$data[80,4] = 2705
$data[80,5] = 0
$data[80,6] = 2
$data[80,7] = 119
$data[81,0] = 'Gabon'
$data[81,1] = 4340
$data[81,2] = 0
$data[81,3] = 1657
$data[81,4] = 2651
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 32
$data[82,0] = 'El Salvador'
$data[82,1] = 4329
$data[82,2] = 129
$data[82,3] = 2310
$data[82,4] = 1933
$data[82,5] = 0
$data[82,6] = 4
$data[82,7] = 86
$data[83,0] = 'Luxemburgo'
$data[83,1] = 4099
$data[83,2] = 8
$data[83,3] = 3944
$data[83,4] = 45
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 110
$data[84,0] = 'Hungria'
$data[84,1] = 4081
$data[84,2] = 2
$data[84,3] = 2581
$data[84,4] = 932
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 568
$data[85,0] = 'Etiopia'
$data[85,1] = 4070
$data[85,2] = 116
$data[85,3] = 1027
$data[85,4] = 2971
$data[85,5] = 0
$data[85,6] = 7
$data[85,7] = 72
$data[86,0] = 'Bulgaria'
$data[86,1] = 3674
$data[86,2] = 0
$data[86,3] = 1941
$data[86,4] = 1543
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 190
$data[87,0] = 'Venezuela'
$data[87,1] = 3483
$data[87,2] = 0
$data[87,3] = 835
$data[87,4] = 2620
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 28
$data[88,0] = 'Bosnia y Herzegovina'
$data[88,1] = 3273
$data[88,2] = 99
$data[88,3] = 2241
$data[88,4] = 863
$data[88,5] = 0
$data[88,6] = 1
$data[88,7] = 169
$data[89,0] = 'Grecia'
$data[89,1] = 3237
$data[89,2] = 10
$data[89,3] = 1374
$data[89,4] = 1674
$data[89,5] = 0
$data[89,6] = 1
$data[89,7] = 189
$data[90,0] = 'Tailandia'
$data[90,1] = 3146
$data[90,2] = 5
$data[90,3] = 3008
$data[90,4] = 80
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 58
$data[91,0] = 'Kirguistan'
$data[91,1] = 2789
$data[91,2] = 132
$data[91,3] = 1961
$data[91,4] = 796
$data[91,5] = 0
$data[91,6] = 1
$data[91,7] = 32
$data[92,0] = 'Somalia'
$data[92,1] = 2719
$data[92,2] = 0
$data[92,3] = 724
$data[92,4] = 1907
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 88
$data[93,0] = 'Mauritania'
$data[93,1] = 2621
$data[93,2] = 197
$data[93,3] = 653
$data[93,4] = 1866
$data[93,5] = 0
$data[93,6] = 5
$data[93,7] = 102
$data[94,0] = 'Republica de Africa Central'
$data[94,1] = 2605
$data[94,2] = 0
$data[94,3] = 417
$data[94,4] = 2169
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 19
$data[95,0] = 'Mayotte'
$data[95,1] = 2394
$data[95,2] = 11
$data[95,3] = 2066
$data[95,4] = 299
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 29
$data[96,0] = 'Cuba'
$data[96,1] = 2305
$data[96,2] = 10
$data[96,3] = 2037
$data[96,4] = 183
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 85
$data[97,0] = 'Croacia'
$data[97,1] = 2280
$data[97,2] = 11
$data[97,3] = 2142
$data[97,4] = 31
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 107
$data[98,0] = 'Maldivas'
$data[98,1] = 2150
$data[98,2] = 13
$data[98,3] = 1769
$data[98,4] = 373
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 8
$data[99,0] = 'Estonia'
$data[99,1] = 1979
$data[99,2] = 2
$data[99,3] = 1755
$data[99,4] = 155
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 69
$data[100,0] = 'Guayana Francesa'
$data[100,1] = 1969
$data[100,2] = 211
$data[100,3] = 840
$data[100,4] = 1124
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 5
$data[101,0] = 'Sri Lanka'
$data[101,1] = 1950
$data[101,2] = 4
$data[101,3] = 1446
$data[101,4] = 493
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 11
$data[102,0] = 'Costa Rica'
$data[102,1] = 1939
$data[102,2] = 0
$data[102,3] = 937
$data[102,4] = 990
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 12
$data[103,0] = 'Mali'
$data[103,1] = 1923
$data[103,2] = 17
$data[103,3] = 1217
$data[103,4] = 598
$data[103,5] = 0
$data[103,6] = 1
$data[103,7] = 108
$data[104,0] = 'Sudan del Sur'
$data[104,1] = 1864
$data[104,2] = 34
$data[104,3] = 122
$data[104,4] = 1708
$data[104,5] = 0
$data[104,6] = 2
$data[104,7] = 34
$data[105,0] = 'Albania'
$data[105,1] = 1838
$data[105,2] = 50
$data[105,3] = 1114
$data[105,4] = 682
$data[105,5] = 0
$data[105,6] = 3
$data[105,7] = 42
$data[106,0] = 'Nicaragua'
$data[106,1] = 1823
$data[106,2] = 0
$data[106,3] = 1238
$data[106,4] = 521
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 64
$data[107,0] = 'Islandia'
$data[107,1] = 1819
$data[107,2] = 3
$data[107,3] = 1801
$data[107,4] = 8
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 10
$data[108,0] = 'Lituania'
$data[108,1] = 1792
$data[108,2] = 8
$data[108,3] = 1462
$data[108,4] = 254
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 76
$data[109,0] = 'Guinea Ecuatorial'
$data[109,1] = 1664
$data[109,2] = 0
$data[109,3] = 515
$data[109,4] = 1117
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 32
$data[110,0] = 'Eslovaquia'
$data[110,1] = 1576
$data[110,2] = 14
$data[110,3] = 1447
$data[110,4] = 101
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 28
$data[111,0] = 'Eslovenia'
$data[111,1] = 1513
$data[111,2] = 2
$data[111,3] = 1359
$data[111,4] = 45
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 109
$data[112,0] = 'Libano'
$data[112,1] = 1510
$data[112,2] = 15
$data[112,3] = 960
$data[112,4] = 518
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 32
$data[113,0] = 'Nueva Zelanda'
$data[113,1] = 1507
$data[113,2] = 0
$data[113,3] = 1482
$data[113,4] = 3
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 22
$data[114,0] = 'Guinea-Bisau'
$data[114,1] = 1492
$data[114,2] = 0
$data[114,3] = 153
$data[114,4] = 1324
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 15
$data[115,0] = 'Madagascar'
$data[115,1] = 1443
$data[115,2] = 40
$data[115,3] = 498
$data[115,4] = 932
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 13
$data[116,0] = 'Zambia'
$data[116,1] = 1416
$data[116,2] = 0
$data[116,3] = 1144
$data[116,4] = 261
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 11
$data[117,0] = 'Paraguay'
$data[117,1] = 1330
$data[117,2] = 0
$data[117,3] = 717
$data[117,4] = 600
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 13
$data[118,0] = 'Sierra Leona'
$data[118,1] = 1298
$data[118,2] = 26
$data[118,3] = 732
$data[118,4] = 513
$data[118,5] = 0
$data[118,6] = 2
$data[118,7] = 53
$data[119,0] = 'Tunez'
$data[119,1] = 1146
$data[119,2] = 14
$data[119,3] = 1014
$data[119,4] = 82
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 50
$data[120,0] = 'Hong Kong'
$data[120,1] = 1128
$data[120,2] = 3
$data[120,3] = 1074
$data[120,4] = 50
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 4
$data[121,0] = 'Letonia'
$data[121,1] = 1110
$data[121,2] = 2
$data[121,3] = 903
$data[121,4] = 177
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 30
$data[122,0] = 'Niger'
$data[122,1] = 1020
$data[122,2] = 0
$data[122,3] = 901
$data[122,4] = 52
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 67
$data[123,0] = 'Jordania'
$data[123,1] = 1008
$data[123,2] = 7
$data[123,3] = 708
$data[123,4] = 291
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 9
$data[124,0] = 'Republica de Chipre'
$data[124,1] = 985
$data[124,2] = 0
$data[124,3] = 818
$data[124,4] = 148
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 19
$data[125,0] = 'Yemen'
$data[125,1] = 909
$data[125,2] = 0
$data[125,3] = 273
$data[125,4] = 388
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 248
$data[126,0] = 'Burkina Faso'
$data[126,1] = 900
$data[126,2] = 1
$data[126,3] = 810
$data[126,4] = 37
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 53
$data[127,0] = 'Georgia'
$data[127,1] = 896
$data[127,2] = 3
$data[127,3] = 741
$data[127,4] = 141
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 14
$data[128,0] = 'Congo'
$data[128,1] = 883
$data[128,2] = 0
$data[128,3] = 391
$data[128,4] = 465
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 27
$data[129,0] = 'Republica del Chad'
$data[129,1] = 858
$data[129,2] = 4
$data[129,3] = 742
$data[129,4] = 42
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 74
$data[130,0] = 'Principado de Andorra'
$data[130,1] = 855
$data[130,2] = 0
$data[130,3] = 792
$data[130,4] = 11
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 52
$data[131,0] = 'Uruguay'
$data[131,1] = 850
$data[131,2] = 0
$data[131,3] = 814
$data[131,4] = 12
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 24
$data[132,0] = 'Cabo Verde'
$data[132,1] = 849
$data[132,2] = 26
$data[132,3] = 377
$data[132,4] = 464
$data[132,5] = 0
$data[132,6] = 1
$data[132,7] = 8
$data[133,0] = 'Uganda'
$data[133,1] = 755
$data[133,2] = 14
$data[133,3] = 492
$data[133,4] = 263
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = 'Crucero'
$data[134,1] = 712
$data[134,2] = 0
$data[134,3] = 651
$data[134,4] = 48
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 13
$data[135,0] = 'San Marino'
$data[135,1] = 696
$data[135,2] = 0
$data[135,3] = 610
$data[135,4] = 44
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 42
$data[136,0] = 'Santo Tome y Principe'
$data[136,1] = 693
$data[136,2] = 5
$data[136,3] = 199
$data[136,4] = 482
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 12
$data[137,0] = 'Estado de Palestina'
$data[137,1] = 673
$data[137,2] = 73
$data[137,3] = 437
$data[137,4] = 233
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 3
$data[138,0] = 'Mozambique'
$data[138,1] = 668
$data[138,2] = 6
$data[138,3] = 177
$data[138,4] = 487
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 4
$data[139,0] = 'Malta'
$data[139,1] = 663
$data[139,2] = 0
$data[139,3] = 613
$data[139,4] = 41
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 9
$data[140,0] = 'Benin'
$data[140,1] = 650
$data[140,2] = 53
$data[140,3] = 247
$data[140,4] = 392
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 11
$data[141,0] = 'Ruanda'
$data[141,1] = 646
$data[141,2] = 0
$data[141,3] = 350
$data[141,4] = 294
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 2
$data[142,0] = 'Jamaica'
$data[142,1] = 638
$data[142,2] = 12
$data[142,3] = 458
$data[142,4] = 170
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 10
$data[143,0] = 'Malaui'
$data[143,1] = 592
$data[143,2] = 0
$data[143,3] = 74
$data[143,4] = 510
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 8
$data[144,0] = 'Suazilandia'
$data[144,1] = 586
$data[144,2] = 0
$data[144,3] = 267
$data[144,4] = 315
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 4
$data[145,0] = 'Liberia'
$data[145,1] = 581
$data[145,2] = 39
$data[145,3] = 250
$data[145,4] = 298
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 33
$data[146,0] = 'Togo'
$data[146,1] = 547
$data[146,2] = 0
$data[146,3] = 353
$data[146,4] = 181
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 13
$data[147,0] = 'Libia'
$data[147,1] = 510
$data[147,2] = 0
$data[147,3] = 81
$data[147,4] = 419
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 10
$data[148,0] = 'Tanzania'
$data[148,1] = 509
$data[148,2] = 0
$data[148,3] = 183
$data[148,4] = 305
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 21
$data[149,0] = 'Reunion'
$data[149,1] = 502
$data[149,2] = 0
$data[149,3] = 460
$data[149,4] = 41
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 1
$data[150,0] = 'Zimbabue'
$data[150,1] = 479
$data[150,2] = 16
$data[150,3] = 63
$data[150,4] = 412
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 4
$data[151,0] = 'Taiwan'
$data[151,1] = 446
$data[151,2] = 0
$data[151,3] = 434
$data[151,4] = 5
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 7
$data[152,0] = 'Vietnam'
$data[152,1] = 349
$data[152,2] = 7
$data[152,3] = 326
$data[152,4] = 23
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = 'Montenegro'
$data[153,1] = 337
$data[153,2] = 0
$data[153,3] = 315
$data[153,4] = 13
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 9
$data[154,0] = 'Mauricio'
$data[154,1] = 337
$data[154,2] = 0
$data[154,3] = 325
$data[154,4] = 2
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 10
$data[155,0] = 'Isla de Man'
$data[155,1] = 336
$data[155,2] = 0
$data[155,3] = 312
$data[155,4] = 0
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 24
$data[156,0] = 'Birmania'
$data[156,1] = 286
$data[156,2] = 23
$data[156,3] = 187
$data[156,4] = 93
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 6
$data[157,0] = 'Surinam'
$data[157,1] = 277
$data[157,2] = 0
$data[157,3] = 74
$data[157,4] = 195
$data[157,5] = 0
$data[157,6] = 1
$data[157,7] = 8
$data[158,0] = 'Martinica'
$data[158,1] = 221
$data[158,2] = 0
$data[158,3] = 98
$data[158,4] = 109
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 14
$data[159,0] = 'Comoras'
$data[159,1] = 210
$data[159,2] = 0
$data[159,3] = 129
$data[159,4] = 76
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 5
$data[160,0] = 'Mongolia'
$data[160,1] = 204
$data[160,2] = 3
$data[160,3] = 132
$data[160,4] = 72
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = 'Islas Caimanes'
$data[161,1] = 193
$data[161,2] = 0
$data[161,3] = 141
$data[161,4] = 51
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 1
$data[162,0] = 'Siria'
$data[162,1] = 187
$data[162,2] = 0
$data[162,3] = 78
$data[162,4] = 102
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 7
$data[163,0] = 'Islas Feroe'
$data[163,1] = 187
$data[163,2] = 0
$data[163,3] = 187
$data[163,4] = 0
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Guyana'
$data[164,1] = 183
$data[164,2] = 0
$data[164,3] = 102
$data[164,4] = 69
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 12
$data[165,0] = 'Gibraltar'
$data[165,1] = 176
$data[165,2] = 0
$data[165,3] = 176
$data[165,4] = 0
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = 'Guadalupe'
$data[166,1] = 171
$data[166,2] = 0
$data[166,3] = 157
$data[166,4] = 0
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 14
$data[167,0] = 'Angola'
$data[167,1] = 166
$data[167,2] = 0
$data[167,3] = 64
$data[167,4] = 94
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 8
$data[168,0] = 'Bermudas'
$data[168,1] = 144
$data[168,2] = 0
$data[168,3] = 128
$data[168,4] = 7
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 9
$data[169,0] = 'Eritrea'
$data[169,1] = 142
$data[169,2] = 0
$data[169,3] = 39
$data[169,4] = 103
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Brunei'
$data[170,1] = 141
$data[170,2] = 0
$data[170,3] = 138
$data[170,4] = 0
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 3
$data[171,0] = 'Camboya'
$data[171,1] = 129
$data[171,2] = 0
$data[171,3] = 126
$data[171,4] = 3
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Trinidad yTobago'
$data[172,1] = 123
$data[172,2] = 0
$data[172,3] = 109
$data[172,4] = 6
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 8
$data[173,0] = 'Burundi'
$data[173,1] = 104
$data[173,2] = 0
$data[173,3] = 75
$data[173,4] = 28
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 1
$data[174,0] = 'Bahamas'
$data[174,1] = 104
$data[174,2] = 0
$data[174,3] = 74
$data[174,4] = 19
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 11
$data[175,0] = 'Aruba'
$data[175,1] = 101
$data[175,2] = 0
$data[175,3] = 98
$data[175,4] = 0
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 3
$data[176,0] = 'Monaco'
$data[176,1] = 99
$data[176,2] = 0
$data[176,3] = 94
$data[176,4] = 1
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 4
$data[177,0] = 'Barbados'
$data[177,1] = 97
$data[177,2] = 0
$data[177,3] = 85
$data[177,4] = 5
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 7
$data[178,0] = 'Botsuana'
$data[178,1] = 89
$data[178,2] = 10
$data[178,3] = 25
$data[178,4] = 63
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 1
$data[179,0] = 'Liechtenstein'
$data[179,1] = 82
$data[179,2] = 0
$data[179,3] = 55
$data[179,4] = 26
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 1
$data[180,0] = 'San Martin (Parte Holandesa)'
$data[180,1] = 77
$data[180,2] = 0
$data[180,3] = 62
$data[180,4] = 0
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 15
$data[181,0] = 'Butan'
$data[181,1] = 67
$data[181,2] = 0
$data[181,3] = 25
$data[181,4] = 42
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'Polinesia Francesa'
$data[182,1] = 60
$data[182,2] = 0
$data[182,3] = 60
$data[182,4] = 0
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Namibia'
$data[183,1] = 45
$data[183,2] = 6
$data[183,3] = 19
$data[183,4] = 26
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Macao'
$data[184,1] = 45
$data[184,2] = 0
$data[184,3] = 45
$data[184,4] = 0
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'San Martin (Parte Francesa)'
$data[185,1] = 42
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 3
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 3
$data[186,0] = 'Puerto Rico'
$data[186,1] = 39
$data[186,2] = 0
$data[186,3] = 1
$data[186,4] = 36
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 2
$data[187,0] = 'Gambia'
$data[187,1] = 36
$data[187,2] = 0
$data[187,3] = 24
$data[187,4] = 11
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 1
$data[188,0] = 'Guam'
$data[188,1] = 32
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 31
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 1
$data[189,0] = 'San Vicente y las Granadinas'
$data[189,1] = 29
$data[189,2] = 0
$data[189,3] = 26
$data[189,4] = 3
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = 'Antigua y Barbuda'
$data[190,1] = 26
$data[190,2] = 0
$data[190,3] = 22
$data[190,4] = 1
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 3
$data[191,0] = 'Timor Oriental'
$data[191,1] = 24
$data[191,2] = 0
$data[191,3] = 24
$data[191,4] = 0
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = 'Curazao'
$data[192,1] = 23
$data[192,2] = 0
$data[192,3] = 19
$data[192,4] = 3
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 1
$data[193,0] = 'Granada'
$data[193,1] = 23
$data[193,2] = 0
$data[193,3] = 23
$data[193,4] = 0
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = 'Belice'
$data[194,1] = 22
$data[194,2] = 0
$data[194,3] = 16
$data[194,4] = 4
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 2
$data[195,0] = 'Nueva Caledonia'
$data[195,1] = 21
$data[195,2] = 0
$data[195,3] = 21
$data[195,4] = 0
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$data[196,0] = 'Santa Lucia'
$data[196,1] = 19
$data[196,2] = 0
$data[196,3] = 18
$data[196,4] = 1
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = 'Laos'
$data[197,1] = 19
$data[197,2] = 0
$data[197,3] = 19
$data[197,4] = 0
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Dominica'
$data[198,1] = 18
$data[198,2] = 0
$data[198,3] = 18
$data[198,4] = 0
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = 'Fiyi'
$data[199,1] = 18
$data[199,2] = 0
$data[199,3] = 18
$data[199,4] = 0
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = 'Islas Virgenes de los Estados Unidos'
$data[200,1] = 17
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 17
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'San Cristobal y Nieves'
$data[201,1] = 15
$data[201,2] = 0
$data[201,3] = 15
$data[201,4] = 0
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = 'Islas Malvinas'
$data[202,1] = 13
$data[202,2] = 0
$data[202,3] = 13
$data[202,4] = 0
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = 'Groenlandia'
$data[203,1] = 13
$data[203,2] = 0
$data[203,3] = 13
$data[203,4] = 0
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'Islas Turcas y Caicos'
$data[204,1] = 12
$data[204,2] = 0
$data[204,3] = 11
$data[204,4] = 0
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 1
$data[205,0] = 'Santa Sede'
$data[205,1] = 12
$data[205,2] = 0
$data[205,3] = 12
$data[205,4] = 0
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = 'Seychelles'
$data[206,1] = 11
$data[206,2] = 0
$data[206,3] = 11
$data[206,4] = 0
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = 'Montserrat'
$data[207,1] = 11
$data[207,2] = 0
$data[207,3] = 10
$data[207,4] = 0
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 1
$data[208,0] = 'Sahara Occidental'
$data[208,1] = 9
$data[208,2] = 0
$data[208,3] = 8
$data[208,4] = 0
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 1
$data[209,0] = 'Papua Nueva Guinea'
$data[209,1] = 8
$data[209,2] = 0
$data[209,3] = 8
$data[209,4] = 0
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = 'Islas Virgenes Britanicas'
$data[210,1] = 8
$data[210,2] = 0
$data[210,3] = 7
$data[210,4] = 0
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 1
$data[211,0] = 'Bonaire, San Eustaquio y Saba'
$data[211,1] = 7
$data[211,2] = 0
$data[211,3] = 7
$data[211,4] = 0
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0
$data[212,0] = 'San Bartolome'
$data[212,1] = 6
$data[212,2] = 0
$data[212,3] = 6
$data[212,4] = 0
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0
$data[213,0] = 'Lesoto'
$data[213,1] = 4
$data[213,2] = 0
$data[213,3] = 2
$data[213,4] = 2
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 0
$data[214,0] = 'Anguila'
$data[214,1] = 3
$data[214,2] = 0
$data[214,3] = 3
$data[214,4] = 0
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 0
$data[215,0] = 'San Pedro y Miquelon'
$data[215,1] = 1
$data[215,2] = 0
$data[215,3] = 1
$data[215,4] = 0
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 0

$ws.Range("A4:H219").Value = $data
